$wb = $excel.ActiveWorkbook

# --- "Programs to include": drop the "Treatment of MAM" row (row 28 in the
#     alphabetically sorted program list) and re-apply the sort so the
#     cached sortState range shrinks from A2:B37 to A2:B36.
$ws1 = $wb.Worksheets.Item("Programs to include")
$ws1.Rows(28).Delete()
$ws1.Sort.SortFields.Clear()
$ws1.Sort.SortFields.Add($ws1.Range("A2:A36"))
$ws1.Sort.SetRange($ws1.Range("A2:B36"))
$ws1.Sort.Header = 0
$ws1.Sort.Apply()
$ws1.Activate()
$ws1.Range("A28").Select() | Out-Null

# --- "Coverage scenario": same program row removed.
$ws2 = $wb.Worksheets.Item("Coverage scenario")
$ws2.Rows(28).Delete()
$ws2.Sort.SortFields.Clear()
$ws2.Sort.SortFields.Add($ws2.Range("A2:A36"))
$ws2.Sort.SetRange($ws2.Range("A2:P36"))
$ws2.Sort.Header = 0
$ws2.Sort.Apply()
$ws2.Activate()
$ws2.Range("A28").Select() | Out-Null

# --- "Budget scenario": same program row removed (no sortState to repair
#     here - this sheet was never saved with a cached sort state).
$ws3 = $wb.Worksheets.Item("Budget scenario")
$ws3.Rows(28).Delete()
$ws3.Activate()
$ws3.Range("A22").Select() | Out-Null

# --- "Optimisation options": the "filter programs" flag moves from the
#     text "0,1,2" to the numeric value 1.
$ws4 = $wb.Worksheets.Item("Optimisation options")
$ws4.Range("F2").Value = 1
$ws4.Activate()
$ws4.Range("F3").Select() | Out-Null

# Restore the originally active sheet/tab ("Budget scenario").
$ws3.Activate()
